# Append the latest daily allocation row (2025-11-23) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 83

$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "11/23/2025"
$ws.Cells.Item($row, 2).Value = 0.2078408625589678
$ws.Cells.Item($row, 3).Value = 0.7921591374410322
